$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 194, pushing existing rows 194-305 down to 195-306.
$ws.Rows(194).Insert()

# Populate the newly inserted row 194 with the new record's data.
$ws.Cells.Item(194, 1).Value = 5
$ws.Cells.Item(194, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(194, 3).Value = "Maule"
$ws.Cells.Item(194, 4).Value = 44719
$ws.Cells.Item(194, 5).Value = 7
$ws.Cells.Item(194, 6).Value = 100112003
$ws.Cells.Item(194, 7).Value = "Ajo"
$ws.Cells.Item(194, 8).Value = "Chino"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 300
$ws.Cells.Item(194, 11).Value = 18000
$ws.Cells.Item(194, 12).Value = 18000
$ws.Cells.Item(194, 13).Value = 18000
$ws.Cells.Item(194, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(194, 15).Value = "China"
$ws.Cells.Item(194, 16).Value = 1800
$ws.Cells.Item(194, 17).Value = 10
$ws.Cells.Item(194, 18).Value = "Hortaliza"
